$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B59").Value = 21313
$ws.Range("C59").Value = 21313
$ws.Range("D59").Value = 15314
$ws.Range("E59").Value = 5999
$ws.Range("G59").Value = 6383
$ws.Range("H59").Value = 8886
$ws.Range("I59").Value = 1934
$ws.Range("J59").Value = 923
$ws.Range("K59").Value = 209
